$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> M1
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl10"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 15.6308825
$ws.Range("H2").Value = 31.261765
$ws.Range("I2").Value = 0.08061597743527853
$ws.Range("J2").Value = 0.05642173194834236
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.7166990000000001
$ws.Range("N2").Value = 2.150097
$ws.Range("O2").Value = 0.2276207788704612
$ws.Range("P2").Value = 0.2276207788704611
$ws.Range("Q2").Value = 11.2026378568675
$ws.Range("R2").Value = 67.21582714120501
$ws.Range("S2").Value = 0.01834987157322162
$ws.Range("T2").Value = 0.01284275857130207

# Row 3: ECs -> M2
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl10"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.6308825
$ws.Range("H3").Value = 31.261765
$ws.Range("I3").Value = 0.08061597743527853
$ws.Range("J3").Value = 0.05642173194834236
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.431954666666666
$ws.Range("N3").Value = 7.295864
$ws.Range("O3").Value = 0.7723792211295388
$ws.Range("P3").Value = 0.7723792211295388
$ws.Range("Q3").Value = 38.01359763999333
$ws.Range("R3").Value = 228.08158583996
$ws.Range("S3").Value = 0.06226610586205691
$ws.Range("T3").Value = 0.04357897337704029

# Row 4: FAPs -> M1
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl10"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 25.268178
$ws.Range("H4").Value = 75.80453399999999
$ws.Range("I4").Value = 0.130320144590595
$ws.Range("J4").Value = 0.1368132316846795
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.7166990000000001
$ws.Range("N4").Value = 2.150097
$ws.Range("O4").Value = 0.2276207788704612
$ws.Range("P4").Value = 0.2276207788704611
$ws.Range("Q4").Value = 18.109677904422
$ws.Range("R4").Value = 162.987101139798
$ws.Range("S4").Value = 0.02966357281422235
$ws.Range("T4").Value = 0.0311415343558516

# Row 5: FAPs -> M2
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl10"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25.268178
$ws.Range("H5").Value = 75.80453399999999
$ws.Range("I5").Value = 0.130320144590595
$ws.Range("J5").Value = 0.1368132316846795
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.431954666666666
$ws.Range("N5").Value = 7.295864
$ws.Range("O5").Value = 0.7723792211295388
$ws.Range("P5").Value = 0.7723792211295388
$ws.Range("Q5").Value = 61.45106340526399
$ws.Range("R5").Value = 553.059570647376
$ws.Range("S5").Value = 0.1006565717763726
$ws.Range("T5").Value = 0.1056716973288279

# Row 6: M1 -> M1
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Cxcl10"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.99944133333333
$ws.Range("H6").Value = 116.998324
$ws.Range("I6").Value = 0.2011388725183283
$ws.Range("J6").Value = 0.2111604407215431
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.7166990000000001
$ws.Range("N6").Value = 2.150097
$ws.Range("O6").Value = 0.2276207788704612
$ws.Range("P6").Value = 0.2276207788704611
$ws.Range("Q6").Value = 27.95086060415867
$ws.Range("R6").Value = 251.557745437428
$ws.Range("S6").Value = 0.04578338682374828
$ws.Range("T6").Value = 0.04806450398366748

# Row 7: M1 -> M2
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Cxcl10"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.99944133333333
$ws.Range("H7").Value = 116.998324
$ws.Range("I7").Value = 0.2011388725183283
$ws.Range("J7").Value = 0.2111604407215431
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.431954666666666
$ws.Range("N7").Value = 7.295864
$ws.Range("O7").Value = 0.7723792211295388
$ws.Range("P7").Value = 0.7723792211295388
$ws.Range("Q7").Value = 94.84487334799287
$ws.Range("R7").Value = 853.603860131936
$ws.Range("S7").Value = 0.15535548569458
$ws.Range("T7").Value = 0.1630959367378756

# Row 8: M2 -> M1
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cxcl10"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 78.37889633333333
$ws.Range("H8").Value = 235.136689
$ws.Range("I8").Value = 0.4042376582518635
$ws.Range("J8").Value = 0.4243784456181134
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.7166990000000001
$ws.Range("N8").Value = 2.150097
$ws.Range("O8").Value = 0.2276207788704612
$ws.Range("P8").Value = 0.2276207788704611
$ws.Range("Q8").Value = 56.17407662320367
$ws.Range("R8").Value = 505.566689608833
$ws.Range("S8").Value = 0.09201289062006048
$ws.Range("T8").Value = 0.09659735232743061

# Row 9: M2 -> M2
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cxcl10"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 78.37889633333333
$ws.Range("H9").Value = 235.136689
$ws.Range("I9").Value = 0.4042376582518635
$ws.Range("J9").Value = 0.4243784456181134
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.431954666666666
$ws.Range("N9").Value = 7.295864
$ws.Range("O9").Value = 0.7723792211295388
$ws.Range("P9").Value = 0.7723792211295388
$ws.Range("Q9").Value = 190.6139227060329
$ws.Range("R9").Value = 1715.525304354296
$ws.Range("S9").Value = 0.3122247676318031
$ws.Range("T9").Value = 0.3277810932906828

# Row 10: Neutro -> M1
$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Cxcl10"
$ws.Range("C10").Value = "Cxcr3"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 23.64038533333333
$ws.Range("H10").Value = 70.921156
$ws.Range("I10").Value = 0.1219248350560686
$ws.Range("J10").Value = 0.1279996332036458
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.7166990000000001
$ws.Range("N10").Value = 2.150097
$ws.Range("O10").Value = 0.2276207788704612
$ws.Range("P10").Value = 0.2276207788704611
$ws.Range("Q10").Value = 16.94304052801467
$ws.Range("R10").Value = 152.487364752132
$ws.Range("S10").Value = 0.02775262591911484
$ws.Range("T10").Value = 0.0291353762049472

# Row 11: Neutro -> M2
$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Cxcl10"
$ws.Range("C11").Value = "Cxcr3"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 23.64038533333333
$ws.Range("H11").Value = 70.921156
$ws.Range("I11").Value = 0.1219248350560686
$ws.Range("J11").Value = 0.1279996332036458
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.431954666666666
$ws.Range("N11").Value = 7.295864
$ws.Range("O11").Value = 0.7723792211295388
$ws.Range("P11").Value = 0.7723792211295388
$ws.Range("Q11").Value = 57.49234543319821
$ws.Range("R11").Value = 517.4311088987839
$ws.Range("S11").Value = 0.09417220913695376
$ws.Range("T11").Value = 0.0988642569986986

# Row 12: sCs -> M1
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Cxcl10"
$ws.Range("C12").Value = "Cxcr3"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 11.9753255
$ws.Range("H12").Value = 23.950651
$ws.Range("I12").Value = 0.06176251214786597
$ws.Range("J12").Value = 0.04322651682367576
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.7166990000000001
$ws.Range("N12").Value = 2.150097
$ws.Range("O12").Value = 0.2276207788704612
$ws.Range("P12").Value = 0.2276207788704611
$ws.Range("Q12").Value = 8.582703810524501
$ws.Range("R12").Value = 51.496222863147
$ws.Range("S12").Value = 0.01405843112009357
$ws.Range("T12").Value = 0.009839253427262168

# Row 13: sCs -> M2
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Cxcl10"
$ws.Range("C13").Value = "Cxcr3"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 11.9753255
$ws.Range("H13").Value = 23.950651
$ws.Range("I13").Value = 0.06176251214786597
$ws.Range("J13").Value = 0.04322651682367576
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.431954666666666
$ws.Range("N13").Value = 7.295864
$ws.Range("O13").Value = 0.7723792211295388
$ws.Range("P13").Value = 0.7723792211295388
$ws.Range("Q13").Value = 29.12344873457733
$ws.Range("R13").Value = 174.740692407464
$ws.Range("S13").Value = 0.0477040810277724
$ws.Range("T13").Value = 0.03338726339641359
